$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4099.222
$ws.Range("I19").Value = 4430.5
$ws.Range("K19").Value = 4430.5
$ws.Range("M19").Value = -4255.5
$ws.Range("H41").Value = 1524.8889
$ws.Range("I41").Value = 441.91666
$ws.Range("J41").Value = 3690.8333
$ws.Range("K41").Value = 441.91666
$ws.Range("L41").Value = 3690.8333
$ws.Range("M41").Value = -1.916659999999979
$ws.Range("N41").Value = -4570.8333
$ws.Range("H53").Value = 2571.3333
$ws.Range("I53").Value = 175.83333
$ws.Range("J53").Value = 4966.8335
$ws.Range("K53").Value = 175.83333
$ws.Range("L53").Value = 4966.8335
$ws.Range("M53").Value = 461.16667
$ws.Range("N53").Value = -6240.8335
$ws.Range("H92").Value = 4381.8486
$ws.Range("I92").Value = 3340.6296
$ws.Range("K92").Value = 3340.6296
$ws.Range("M92").Value = -2092.6296
$ws.Range("H98").Value = 7199
$ws.Range("I98").Value = 1048.9166
$ws.Range("J98").Value = 25649.25
$ws.Range("K98").Value = 1048.9166
$ws.Range("L98").Value = 25649.25
$ws.Range("M98").Value = 449.0834
$ws.Range("N98").Value = -28645.25
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()  # was -21506
$ws.Range("N105").ClearContents()  # was -41651.332
$ws.Range("H120").Value = 50560.75
$ws.Range("J120").Value = 50560.75
$ws.Range("L120").Value = 50560.75
$ws.Range("N120").Value = -60236.75
$ws.Range("H122").Value = 7199
$ws.Range("I122").Value = 1048.9166
$ws.Range("J122").Value = 25649.25
$ws.Range("K122").Value = 3146.7498
$ws.Range("L122").Value = 76947.75
$ws.Range("M122").Value = -696.7498000000001
$ws.Range("N122").Value = -81847.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19665.834
$ws.Range("I2").Value = 3998.75
$ws.Range("K2").Value = 3998.75
$ws.Range("M2").Value = -3885.75
$ws.Range("H32").Value = 1212.0822
$ws.Range("I32").Value = 1125.254
$ws.Range("K32").Value = 1125.254
$ws.Range("M32").Value = -838.2539999999999
$ws.Range("H61").Value = 7112.3184
$ws.Range("I61").Value = 6553.9443
$ws.Range("J61").Value = 9625
$ws.Range("K61").Value = 6553.9443
$ws.Range("L61").Value = 9625
$ws.Range("M61").Value = -6341.9443
$ws.Range("N61").Value = -10049
$ws.Range("H74").Value = 9262563
$ws.Range("I74").Value = 11113109
$ws.Range("J74").Value = 9833.666999999999
$ws.Range("K74").Value = 11113109
$ws.Range("L74").Value = 9833.666999999999
$ws.Range("M74").Value = -11112235
$ws.Range("N74").Value = -11581.667
$ws.Range("H77").Value = 9262563
$ws.Range("I77").Value = 11113109
$ws.Range("J77").Value = 9833.666999999999
$ws.Range("K77").Value = 55565545
$ws.Range("L77").Value = 49168.335
$ws.Range("M77").Value = -55561177
$ws.Range("N77").Value = -57904.335
$ws.Range("H116").Value = 19665.834
$ws.Range("I116").Value = 3998.75
$ws.Range("K116").Value = 3998.75
$ws.Range("M116").Value = -1704.75
$ws.Range("H132").Value = 7062.1523
$ws.Range("I132").Value = 6237.278
$ws.Range("J132").Value = 10031.7
$ws.Range("K132").Value = 18711.834
$ws.Range("L132").Value = 30095.1
$ws.Range("M132").Value = -16181.834
$ws.Range("N132").Value = -35155.10000000001
$ws.Range("H136").Value = 7112.3184
$ws.Range("I136").Value = 6553.9443
$ws.Range("J136").Value = 9625
$ws.Range("K136").Value = 19661.8329
$ws.Range("L136").Value = 28875
$ws.Range("M136").Value = -17111.8329
$ws.Range("N136").Value = -33975

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19665.834
$ws.Range("I3").Value = 3998.75
$ws.Range("K3").Value = 3998.75
$ws.Range("M3").Value = -3884.75
$ws.Range("H25").Value = 2750
$ws.Range("I25").Value = 2750
$ws.Range("K25").Value = 2750
$ws.Range("M25").Value = -2515

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2963.7727
$ws.Range("I105").Value = 2298.3
$ws.Range("K105").Value = 2298.3
$ws.Range("M105").Value = -551.3000000000002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5817.067
$ws.Range("I68").Value = 998.75
$ws.Range("J68").Value = 6558.346
$ws.Range("K68").Value = 2996.25
$ws.Range("L68").Value = 19675.038
$ws.Range("M68").Value = -2185.25
$ws.Range("N68").Value = -21297.038
$ws.Range("H71").Value = 5817.067
$ws.Range("I71").Value = 998.75
$ws.Range("J71").Value = 6558.346
$ws.Range("K71").Value = 8988.75
$ws.Range("L71").Value = 59025.11399999999
$ws.Range("M71").Value = -4932.75
$ws.Range("N71").Value = -67137.114
$ws.Range("H109").Value = 4631.375
$ws.Range("I109").Value = 2210.2
$ws.Range("K109").Value = 6630.599999999999
$ws.Range("M109").Value = -5590.599999999999
$ws.Range("H114").Value = 83339256
$ws.Range("J114").Value = 125007500
$ws.Range("L114").Value = 375022500
$ws.Range("N114").Value = -375029008
$ws.Range("H119").Value = 22576.334
$ws.Range("I119").Value = 23285.2
$ws.Range("K119").Value = 69855.60000000001
$ws.Range("M119").Value = -65017.60000000001
$ws.Range("H131").Value = 5721208
$ws.Range("I131").Value = 1732.15
$ws.Range("J131").Value = 13891888
$ws.Range("K131").Value = 5196.450000000001
$ws.Range("L131").Value = 41675664
$ws.Range("M131").Value = -156.4500000000007
$ws.Range("N131").Value = -41685744
$ws.Range("H132").Value = 4356.2144
$ws.Range("J132").Value = 4563.353
$ws.Range("L132").Value = 41070.177
$ws.Range("N132").Value = -46130.177

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5389.2607
$ws.Range("I80").Value = 3392.4
$ws.Range("K80").Value = 3392.4
$ws.Range("M80").Value = -2394.4
$ws.Range("H83").Value = 5389.2607
$ws.Range("I83").Value = 3392.4
$ws.Range("K83").Value = 16962
$ws.Range("M83").Value = -11970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2349.4707
$ws.Range("I61").Value = 2003
$ws.Range("K61").Value = 2003
$ws.Range("M61").Value = -1801
$ws.Range("H113").Value = 2349.4707
$ws.Range("I113").Value = 2003
$ws.Range("K113").Value = 2003
$ws.Range("M113").Value = 167

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()  # was -10460
$ws.Range("H30").Value = 303.8
$ws.Range("J30").Value = 297.25
$ws.Range("L30").Value = 297.25
$ws.Range("N30").Value = -511.25
$ws.Range("H81").Value = 7740.4
$ws.Range("I81").Value = 1849.75
$ws.Range("J81").Value = 11667.5
$ws.Range("K81").Value = 3699.5
$ws.Range("L81").Value = 23335
$ws.Range("M81").Value = -2638.5
$ws.Range("N81").Value = -25457
$ws.Range("H84").Value = 7740.4
$ws.Range("I84").Value = 1849.75
$ws.Range("J84").Value = 11667.5
$ws.Range("K84").Value = 18497.5
$ws.Range("L84").Value = 116675
$ws.Range("M84").Value = -13193.5
$ws.Range("N84").Value = -127283
$ws.Range("H136").Value = 2730.7715
$ws.Range("I136").Value = 1562.96
$ws.Range("J136").Value = 5650.3
$ws.Range("K136").Value = 4688.88
$ws.Range("L136").Value = 16950.9
$ws.Range("M136").Value = -2138.88
$ws.Range("N136").Value = -22050.9
